# "perbaikan final untuk app"
# The import template used to store a single "tanggal_meter" (meter date) column.
# The app now wants the date split into an Indonesian month name ("bulan_meter")
# and a year ("tahun_meter") instead, so rebuild the sheet1 header/data layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$indoMonths = @("Januari","Februari","Maret","April","Mei","Juni","Juli", `
                "Agustus","September","Oktober","November","Desember")

# Remember the existing "tanggal_meter" date values (column A) before we touch
# the layout, so we can translate each one into month name + year.
$date2 = $ws.Range("A2").Value()
$date3 = $ws.Range("A3").Value()

$bulan2 = $indoMonths[$date2.Month - 1]
$tahun2 = $date2.Year
$bulan3 = $indoMonths[$date3.Month - 1]
$tahun3 = $date3.Year

# Insert a new column before column A; no_sambung/tanggal_meter/meter all shift
# one column to the right (tanggal_meter ends up in the new column B).
$ws.Range("A1").EntireColumn.Insert() | Out-Null

# New header row: bulan_meter | tahun_meter | meter | no_sambung
$ws.Range("A1").Value = "bulan_meter"
$ws.Range("B1").Value = "tahun_meter"

# Replace the old date values with the month name (col A) and year (col B).
$ws.Range("A2").Value = $bulan2
$ws.Range("B2").Value = $tahun2
$ws.Range("A3").Value = $bulan3
$ws.Range("B3").Value = $tahun3

# The old column used a custom yyyy-mm-dd date format; the new month/year
# columns should just show as plain General values, not dates.
$ws.Range("A2:B3").NumberFormat = "General"

$ws.Range("B2").Select() | Out-Null

Write-Host "Rebuilt sheet1 with bulan_meter/tahun_meter columns."
